$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.914.34'
$ws.Range("E2").Value = '  -1.30%  '
$ws.Range("D3").Value = '1.636.80'
$ws.Range("E3").Value = '  -0.54%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").Value = "'215.42"
$ws.Range("E5").Value = '  -0.68%  '
$ws.Range("E6").Value = '  +0.26%  '
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("E8").Value = '  -1.08%  '
$ws.Range("E9").Value = '  -0.13%  '
$ws.Range("D10").Value = "'19.57"
$ws.Range("E10").Value = '  -1.75%  '
$ws.Range("E11").Value = '  -0.24%  '
$ws.Range("E12").Value = '  -0.46%  '
$ws.Range("D13").Value = '1.863.79'
$ws.Range("E13").Value = '  -0.48%  '
$ws.Range("D14").Value = '1.660.99'
$ws.Range("E14").Value = '  +1.11%  '
$ws.Range("E15").Value = '  -0.65%  '
$ws.Range("D16").Value = '0.0₃0765'
$ws.Range("E16").Value = '  -0.06%  '
$ws.Range("D17").Value = "'62.79"
$ws.Range("E17").Value = '  -0.89%  '
$ws.Range("D18").Value = '25.948.45'
$ws.Range("E18").Value = '  -1.14%  '
$ws.Range("D20").Value = "'192.80"
$ws.Range("E20").Value = '  -1.44%  '
$ws.Range("E21").Value = '  -2.10%  '
$ws.Range("E22").Value = '  -1.68%  '
$ws.Range("E23").Value = '  -1.07%  '
$ws.Range("E24").Value = '  +4.45%  '
$ws.Range("E25").Value = '  -0.09%  '
$ws.Range("D26").Value = "'143.30"
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("E27").Value = '  +0.08%  '
$ws.Range("E28").Value = '  -0.99%  '
$ws.Range("E29").Value = '  -0.68%  '
$ws.Range("D30").Value = "'1.24"
$ws.Range("E30").Value = '  -0.78%  '
$ws.Range("D31").Value = "'0.0501"
$ws.Range("E31").Value = '  -0.38%  '
$ws.Range("E32").Value = '  -2.11%  '
$ws.Range("E33").Value = '  -0.24%  '
$ws.Range("E34").Value = '  -4.24%  '
$ws.Range("E35").Value = '  +1.57%  '
$ws.Range("E36").Value = '  -1.42%  '
$ws.Range("D37").Value = '1.133.82'
$ws.Range("E37").Value = '  -0.39%  '
$ws.Range("E39").Value = '  -1.32%  '
$ws.Range("E41").Value = '  -0.71%  '
$ws.Range("D42").Value = "'99.27"
$ws.Range("E42").Value = '  -1.30%  '
$ws.Range("E43").Value = '  -0.66%  '
$ws.Range("D44").Value = '1.773.47'
$ws.Range("E45").Value = '  +1.92%  '
$ws.Range("D46").Value = "'56.61"
$ws.Range("E46").Value = '  -0.42%  '
$ws.Range("E47").Value = '  +2.26%  '
$ws.Range("E48").Value = '  -0.74%  '
$ws.Range("D49").Value = "'7.69"
$ws.Range("E49").Value = '  -0.19%  '
$ws.Range("E50").Value = '  -0.89%  '
$ws.Range("D51").Value = "'0.0958"
$ws.Range("E51").Value = '  -1.43%  '

# Restore default (unstyled) cell style for text-forced numeric-looking price cells
$ws.Range("D5").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"
